$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.952.42"
$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("D3").Value = "2.788.73"
$ws.Range("E3").Value = "  -1.00%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "358.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.565"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.47%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.595"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.15"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.05%  "
$ws.Range("E11").Value = "  +0.25%  "
$ws.Range("E12").Value = "  +0.94%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.53"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.13%  "
$ws.Range("E14").Value = "  -1.15%  "
$ws.Range("D15").Value = "3.227.65"
$ws.Range("E15").Value = "  -0.97%  "
$ws.Range("D16").Value = "2.810.83"
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.936"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.29%  "
$ws.Range("D18").Value = "51.906.73"
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.41"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.37%  "
$ws.Range("E22").Value = "  -1.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "274.65"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.01%  "
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.70"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("E28").Value = "  -1.02%  "
$ws.Range("E29").Value = "  +4.95%  "
$ws.Range("E30").Value = "  -1.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0464"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "51.57"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.84%  "
$ws.Range("E34").Value = "  -1.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0846"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.27"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.39%  "
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("E38").Value = "  +1.55%  "
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.98%  "
$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.05"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.18%  "
$ws.Range("E41").Value = "  +3.17%  "
$ws.Range("E42").Value = "  -1.38%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "122.39"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.49%  "
$ws.Range("E44").Value = "  -2.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.96"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.20%  "
$ws.Range("D46").Value = "2.074.86"
$ws.Range("E46").Value = "  -0.20%  "
$ws.Range("E47").Value = "  -1.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.74"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.940"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.96"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.66%  "
